# 11.5.1.xlsx — add a new "2021" data column (R) that mirrors the existing
# "2020" column (Q): same per-row formatting, new figures for the extra year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is the blank, bottom-bordered spacer row above the header — just
# needs the same formatting as Q3, no value.
$ws.Range("Q3").Copy($ws.Range("R3"))

# Header row: new year label.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# Data rows 5-34: copy Q's formatting into R, then overwrite with the 2021
# figures (numbers, or "-" where the source sheet has no data).
$values = @{
    5  = 109
    6  = 74
    7  = 35
    8  = 36
    9  = 35
    10 = 1
    11 = 15
    12 = 8
    13 = 7
    14 = 12
    15 = 7
    16 = 5
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 17
    21 = 8
    22 = 9
    23 = 9
    24 = 7
    25 = 2
    26 = 20
    27 = 9
    28 = 11
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($r in 5..34) {
    $ws.Range("Q$r").Copy($ws.Range("R$r"))
    $ws.Range("R$r").Value = $values[$r]
}

# Move the active selection the way the original edit left it.
$ws.Range("R1").Select()
